# Auto-generated edit script applying the diff's cell-value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value2 = 215
$ws.Range("F4").Value2 = 399
$ws.Range("F5").Value2 = 197
$ws.Range("F6").Value2 = 790
$ws.Range("F7").Value2 = 88
$ws.Range("F8").Value2 = 10120
$ws.Range("F10").Value2 = 3490
$ws.Range("F14").Value2 = 2777
$ws.Range("F16").Value2 = 508
$ws.Range("F17").Value2 = 2148
$ws.Range("F21").Value2 = 382
$ws.Range("F23").Value2 = 134
$ws.Range("I23").Value2 = "//i1.hdslb.com/bfs/openplatform/202408/eUahHqE51723694029087.png"
$ws.Range("F24").Value2 = 310
$ws.Range("F26").Value2 = 213
$ws.Range("F29").Value2 = 7
$ws.Range("F30").Value2 = 1248
$ws.Range("F31").Value2 = 101
$ws.Range("F34").Value2 = 3075
$ws.Range("F35").Value2 = 2964
$ws.Range("G35").Value2 = 90
$ws.Range("F36").Value2 = 21
$ws.Range("F38").Value2 = 1028
$ws.Range("F39").Value2 = 392
$ws.Range("F41").Value2 = 1290
$ws.Range("F42").Value2 = 87
$ws.Range("F43").Value2 = 104
$ws.Range("F44").Value2 = 70

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value2 = 9
$ws.Range("F4").Value2 = 177
$ws.Range("F13").Value2 = 5
$ws.Range("F16").Value2 = 174

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value2 = 982
$ws.Range("F5").Value2 = 1982

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value2 = 982
$ws.Range("F5").Value2 = 399
$ws.Range("F6").Value2 = 9
$ws.Range("F8").Value2 = 197
$ws.Range("F9").Value2 = 88
$ws.Range("F10").Value2 = 10120
$ws.Range("F12").Value2 = 3490
$ws.Range("F16").Value2 = 508
$ws.Range("F17").Value2 = 2148
$ws.Range("F21").Value2 = 134
$ws.Range("I21").Value2 = "//i1.hdslb.com/bfs/openplatform/202408/eUahHqE51723694029087.png"
$ws.Range("F22").Value2 = 310
$ws.Range("F24").Value2 = 213
$ws.Range("F27").Value2 = 7
$ws.Range("F28").Value2 = 1248
$ws.Range("F33").Value2 = 3075
$ws.Range("F34").Value2 = 2964
$ws.Range("G34").Value2 = 90
$ws.Range("F35").Value2 = 21
$ws.Range("F36").Value2 = 1028
$ws.Range("F39").Value2 = 392
$ws.Range("F40").Value2 = 5
$ws.Range("F44").Value2 = 87
$ws.Range("F45").Value2 = 70
$ws.Range("F49").Value2 = 174
